$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$newValues = @(
    "74-26=",
    "17+74=",
    "70+5=",
    "55+37=",
    "47-5=",
    "72+20=",
    "97-38=",
    "84-73=",
    "64-2=",
    "90-18=",
    "92-5=",
    "36+21=",
    "6+11=",
    "16+66=",
    "25+5=",
    "21-8=",
    "57-6=",
    "37-20=",
    "71-67=",
    "55+13=",
    "41+18=",
    "61+35=",
    "84-59=",
    "18+36=",
    "69-56=",
    "17+44=",
    "76-3=",
    "94-14=",
    "12-8=",
    "25+15=",
    "3+60=",
    "73+14=",
    "65-20=",
    "22-14=",
    "50-4=",
    "18+42=",
    "44-8=",
    "48+19=",
    "25-5=",
    "98-12=",
    "25+10=",
    "30+16=",
    "99-4=",
    "93-93=",
    "99-46=",
    "58-16=",
    "44+37=",
    "50+1=",
    "79-51=",
    "48+29=",
    "77-59=",
    "99-66=",
    "63-2=",
    "51+3=",
    "81+13=",
    "27+38=",
    "36-32=",
    "81-47=",
    "97-58=",
    "72-61=",
    "10+89=",
    "87-41=",
    "24+10=",
    "45+8=",
    "74-3=",
    "39+36=",
    "34+14=",
    "95+4=",
    "63+19=",
    "22-4=",
    "83-43=",
    "76-51=",
    "88-12=",
    "40+47=",
    "49+5=",
    "79-6=",
    "12+51=",
    "75-61=",
    "86-76=",
    "5+79=",
    "27+20=",
    "11+18=",
    "6+47=",
    "87-2=",
    "14+58=",
    "19+71=",
    "37-20=",
    "52+30=",
    "11+58=",
    "75-21=",
    "19+39=",
    "68-32=",
    "8+53=",
    "54+23=",
    "42-35=",
    "34+3=",
    "30+20=",
    "56-41=",
    "1+85=",
    "93-10="
)

$idx = 0
foreach ($row in $table.Rows) {
    foreach ($cell in $row.Cells) {
        if ($idx -lt $newValues.Count) {
            $cellRange = $cell.Range
            $cellRange.MoveEnd(1, -1) | Out-Null
            $cellRange.Text = $newValues[$idx]
        }
        $idx = $idx + 1
    }
}

Write-Output "Replaced $idx cells"
